$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "Harvard"
$ws.Range("B10").Value = 0.8
$ws.Range("C10").Value = 0.2

$ws.Range("C10").Select()
